# New crime data collected - update weekly CompStat figures for the 123rd Precinct.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Header text: volume/issue number and the reporting week dates.
# ---------------------------------------------------------------------------
$ws.Range("A8").Value = "Volume 30   Number  17"
$ws.Range("C9").Value = "Report Covering the Week  4/24/2023  Through  4/30/2023"

# ---------------------------------------------------------------------------
# Cells that flip from a numeric value to the "no data" placeholder text
# ("0" / "***.*"). Copying from a cell that already holds that placeholder
# text keeps both the shared-string type AND the original (General format)
# cell style intact.
# ---------------------------------------------------------------------------
$ws.Range("C16").Copy($ws.Range("D16"))   # D16: 1 -> "0"
$ws.Range("E17").Copy($ws.Range("E16"))   # E16: -100 -> "***.*"

# ---------------------------------------------------------------------------
# Cells that flip from the placeholder text back to a real numeric value.
# Setting .Value first, then restoring NumberFormat from a cell that already
# carries the desired numeric style, avoids the engine creating a stray
# "@" (text) number format and keeps the original style id.
# ---------------------------------------------------------------------------
$ws.Range("C17").Value = 3
$ws.Range("C17").NumberFormat = $ws.Range("F17").NumberFormat

$ws.Range("C18").Value = 1
$ws.Range("C18").NumberFormat = $ws.Range("F17").NumberFormat

$ws.Range("F18").Value = 1
$ws.Range("F18").NumberFormat = $ws.Range("F17").NumberFormat

$ws.Range("C20").Value = 3
$ws.Range("C20").NumberFormat = $ws.Range("F17").NumberFormat

$ws.Range("C27").Value = 1
$ws.Range("C27").NumberFormat = $ws.Range("F17").NumberFormat

$ws.Range("D27").Value = 1
$ws.Range("D27").NumberFormat = $ws.Range("F17").NumberFormat

$ws.Range("E27").Value = 0
$ws.Range("E27").NumberFormat = $ws.Range("H16").NumberFormat

$ws.Range("F27").Value = 1
$ws.Range("F27").NumberFormat = $ws.Range("F17").NumberFormat

# ---------------------------------------------------------------------------
# Row 16 - Robbery
# ---------------------------------------------------------------------------
$ws.Range("M16").Value = -33.333333333333

# ---------------------------------------------------------------------------
# Row 17 - Fel. Assault
# ---------------------------------------------------------------------------
$ws.Range("F17").Value = 9
$ws.Range("H17").Value = 350
$ws.Range("I17").Value = 28
$ws.Range("K17").Value = 115.384615384615
$ws.Range("L17").Value = 154.545454545455
$ws.Range("M17").Value = 86.666666666666
$ws.Range("N17").Value = 0

# ---------------------------------------------------------------------------
# Row 18 - Burglary
# ---------------------------------------------------------------------------
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = 0
$ws.Range("G18").Value = 4
$ws.Range("H18").Value = -75
$ws.Range("I18").Value = 20
$ws.Range("J18").Value = 6
$ws.Range("K18").Value = 233.333333333333
$ws.Range("L18").Value = 53.846153846153
$ws.Range("M18").Value = -48.717948717948
$ws.Range("N18").Value = -76.744186046511

# ---------------------------------------------------------------------------
# Row 19 - Gr. Larceny
# ---------------------------------------------------------------------------
$ws.Range("C19").Value = 3
$ws.Range("D19").Value = 4
$ws.Range("E19").Value = -25
$ws.Range("F19").Value = 26
$ws.Range("G19").Value = 15
$ws.Range("H19").Value = 73.333333333333
$ws.Range("I19").Value = 93
$ws.Range("J19").Value = 87
$ws.Range("K19").Value = 6.896551724137
$ws.Range("L19").Value = 86
$ws.Range("M19").Value = 89.795918367346
$ws.Range("N19").Value = 63.157894736842

# ---------------------------------------------------------------------------
# Row 20 - G.L.A.
# ---------------------------------------------------------------------------
$ws.Range("E20").Value = 200
$ws.Range("F20").Value = 7
$ws.Range("G20").Value = 9
$ws.Range("H20").Value = -22.222222222222
$ws.Range("I20").Value = 25
$ws.Range("J20").Value = 35
$ws.Range("K20").Value = -28.571428571428
$ws.Range("L20").Value = 108.333333333333
$ws.Range("M20").Value = 108.333333333333
$ws.Range("N20").Value = -89.224137931034

# ---------------------------------------------------------------------------
# Row 21 - TOTAL
# ---------------------------------------------------------------------------
$ws.Range("C21").Value = 10
$ws.Range("E21").Value = 66.666666666666
$ws.Range("F21").Value = 44
$ws.Range("G21").Value = 32
$ws.Range("H21").Value = 37.5
$ws.Range("I21").Value = 171
$ws.Range("J21").Value = 145
$ws.Range("K21").Value = 17.931034482758
$ws.Range("L21").Value = 94.318181818181
$ws.Range("M21").Value = 39.024390243902
$ws.Range("N21").Value = -59.574468085106

# ---------------------------------------------------------------------------
# Row 24 - Petit Larceny
# ---------------------------------------------------------------------------
$ws.Range("D24").Value = 10
$ws.Range("E24").Value = -60
$ws.Range("F24").Value = 24
$ws.Range("G24").Value = 34
$ws.Range("H24").Value = -29.411764705882
$ws.Range("I24").Value = 160
$ws.Range("J24").Value = 133
$ws.Range("K24").Value = 20.300751879699
$ws.Range("L24").Value = 122.222222222222
$ws.Range("M24").Value = -3.614457831325

# ---------------------------------------------------------------------------
# Row 25 - Misd. Assault
# ---------------------------------------------------------------------------
$ws.Range("C25").Value = 4
$ws.Range("D25").Value = 2
$ws.Range("E25").Value = 100
$ws.Range("G25").Value = 12
$ws.Range("H25").Value = 8.333333333333
$ws.Range("I25").Value = 63
$ws.Range("J25").Value = 61
$ws.Range("K25").Value = 3.278688524590
$ws.Range("L25").Value = 85.294117647058
$ws.Range("M25").Value = 3.278688524590

# ---------------------------------------------------------------------------
# Row 27 - Other Sex Crimes
# ---------------------------------------------------------------------------
$ws.Range("H27").Value = -50
$ws.Range("I27").Value = 4
$ws.Range("J27").Value = 9
$ws.Range("K27").Value = -55.555555555555
$ws.Range("L27").Value = 33.333333333333
